$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has columns:
#  A QUIZ_ID, B CLASS_ID, C QUIZ_NAME, D START_TIME, E TIME_BEGIN,
#  F END_TIME, G TIME_FINISH, H LENGTH, I WEIGHT, J IS_DELETED
#
# TIME_BEGIN (E) and TIME_FINISH (G) are dropped; the time-of-day that used
# to live in TIME_FINISH is folded into the END_TIME (old F, new E) value as
# a combined date+time serial number.
#
# Delete G first, then E, so neither deletion needs to account for the
# other having already shifted column letters.
$ws.Columns("G").Delete()
$ws.Columns("E").Delete()

# After the two deletions, the surviving columns are:
#  A QUIZ_ID, B CLASS_ID, C QUIZ_NAME, D START_TIME, E END_TIME,
#  F LENGTH, G WEIGHT, H IS_DELETED

# Fold the former TIME_FINISH time-of-day into the END_TIME column (E) as a
# combined date+time serial value for every data row.
$ws.Range("E2").Value = 44941.000138888892
$ws.Range("E3").Value = 44941.508472222224
$ws.Range("E4").Value = 44946.25708333333
$ws.Range("E5").Value = 44946.229166666664
$ws.Range("E6").Value = 44958.333333333336
$ws.Range("E7").Value = 44958.375
$ws.Range("E8").Value = 44969.416666666664
$ws.Range("E9").Value = 45052.458333333336

# Match the saved selection state from the edited workbook.
[void]$ws.Range("D2").Select()
